# Commit: Add PF/1.0.4 to meta-sheet
#
# The meta-sheet has a header row (dev2/sit2/uat2/prod) and a row of
# "PF/1.0.0" values. This adds a new row documenting that "PF/1.0.4"
# applies to (is only relevant for / marked "X" against) all four
# environment columns.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A3").Value2 = "PF/1.0.4"
$ws.Range("B3").Value2 = "X"
$ws.Range("C3").Value2 = "X"
$ws.Range("D3").Value2 = "X"

# New row should use the plain default style (no special alignment
# formatting carried over from the column style), matching the rest
# of the sheet's look for this new entry.
$ws.Range("A3:D3").Style = "Normal"
